# set node_rel and node_loc equal to fully reproduce results behaviour
# Fills in previously-zero report cells (B2:N13) on all 4 sheets (CAP, CAP_NEW, INVESTMENT, REMOVAL)
# with their newly computed values.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: CAP ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 0.03192683931804503
$ws.Range("C2").Value = 0.6421116143994228
$ws.Range("D2").Value = 1.266626409862353
$ws.Range("E2").Value = 2.444955656598223
$ws.Range("F2").Value = 4.668217337186446
$ws.Range("G2").Value = 8.863048717349876
$ws.Range("H2").Value = 16.77782117155247
$ws.Range("I2").Value = 31.71134705724933
$ws.Range("J2").Value = 102.3379908040818
$ws.Range("K2").Value = 329.9420322643048
$ws.Range("L2").Value = 1063.427257397211
$ws.Range("M2").Value = 2545.467977962364
$ws.Range("N2").Value = 2032.82116599549
$ws.Range("E3").Value = 0.1951107636166253
$ws.Range("F3").Value = 1.144474502355429
$ws.Range("G3").Value = 2.935725994351259
$ws.Range("H3").Value = 6.315444195764529
$ws.Range("I3").Value = 12.69226798181876
$ws.Range("J3").Value = 42.92921987589357
$ws.Range("K3").Value = 140.3719424438682
$ws.Range("L3").Value = 454.3944740377329
$ws.Range("M3").Value = 1466.375130467265
$ws.Range("N3").Value = 1375.164324299622
$ws.Range("C4").Value = 0.00832013724340157
$ws.Range("D4").Value = 0.09894068077803235
$ws.Range("E4").Value = 0.2699227577083377
$ws.Range("F4").Value = 0.5925302977991638
$ws.Range("G4").Value = 1.201223482441415
$ws.Range("H4").Value = 2.34970064737567
$ws.Range("I4").Value = 4.516637671897544
$ws.Range("J4").Value = 14.77278896785996
$ws.Range("K4").Value = 47.82464250285461
$ws.Range("L4").Value = 154.3387714947697
$ws.Range("M4").Value = 237.1757854849172
$ws.Range("N4").Value = 420.6496069519287
$ws.Range("C5").Value = 0.1705807431751008
$ws.Range("D5").Value = 0.4643803156869322
$ws.Range("E5").Value = 1.018718870168363
$ws.Range("F5").Value = 2.064640177532373
$ws.Range("G5").Value = 4.038075675305758
$ws.Range("H5").Value = 7.761537121642568
$ws.Range("I5").Value = 14.78693276188491
$ws.Range("J5").Value = 48.02499393272711
$ws.Range("K5").Value = 155.1392025577843
$ws.Range("L5").Value = 500.3294238044214
$ws.Range("M5").Value = 1373.86758678948
$ws.Range("N5").Value = 1477.091646810583
$ws.Range("C6").Value = 0.4927437130866399
$ws.Range("D6").Value = 1.17978390228413
$ws.Range("E6").Value = 2.476085534360359
$ws.Range("F6").Value = 4.921936516216769
$ws.Range("G6").Value = 9.536747511915785
$ws.Range("H6").Value = 18.24393396113291
$ws.Range("I6").Value = 34.67257988663013
$ws.Range("J6").Value = 112.389976150087
$ws.Range("K6").Value = 362.8449396442537
$ws.Range("L6").Value = 1169.970393593092
$ws.Range("M6").Value = 2710.197870741217
$ws.Range("N6").Value = 2150.600259812235
$ws.Range("B7").Value = 0.02444303768772006
$ws.Range("C7").Value = 0.7922271641676304
$ws.Range("D7").Value = 1.681085180002876
$ws.Range("E7").Value = 3.358174984541308
$ws.Range("F7").Value = 6.522493877299199
$ws.Range("G7").Value = 12.49290406521413
$ws.Range("H7").Value = 23.75782363243607
$ws.Range("I7").Value = 45.01237882411158
$ws.Range("J7").Value = 145.5449454584486
$ws.Range("K7").Value = 469.5249231332411
$ws.Range("L7").Value = 1513.594814593753
$ws.Range("M7").Value = 4123.44091719734
$ws.Range("N7").Value = 5594.135919941928
$ws.Range("B8").Value = 0.0252984265455309
$ws.Range("C8").Value = 2.096199403634548
$ws.Range("D8").Value = 4.759825782014541
$ws.Range("E8").Value = 9.785533558394343
$ws.Range("F8").Value = 19.2679965864584
$ws.Range("G8").Value = 37.15942781952787
$ws.Range("H8").Value = 70.91682931113232
$ws.Range("I8").Value = 134.6100096174165
$ws.Range("J8").Value = 435.8999765451067
$ws.Range("K8").Value = 1406.848197870359
$ws.Range("L8").Value = 4535.861945350807
$ws.Range("M8").Value = 9405.134321401261
$ws.Range("N8").Value = 9315.136706420724
$ws.Range("D9").Value = 0.04530537974837011
$ws.Range("E9").Value = 0.2929358916521486
$ws.Range("F9").Value = 0.760163052129754
$ws.Range("G9").Value = 1.641723316289733
$ws.Range("H9").Value = 3.305043784542053
$ws.Range("I9").Value = 6.443382809774082
$ws.Range("J9").Value = 21.30624536700153
$ws.Range("K9").Value = 69.20385746108481
$ws.Range("L9").Value = 223.5604785683471
$ws.Range("M9").Value = 720.9958503537005
$ws.Range("N9").Value = 1346.581622336518
$ws.Range("K10").Value = 0.6438220879409279
$ws.Range("L10").Value = 2.718627025890233
$ws.Range("M10").Value = 2.396715981919769
$ws.Range("N10").Value = 1.435002391238348
$ws.Range("K11").Value = 0.4092284375570817
$ws.Range("L11").Value = 1.728023177433443
$ws.Range("M11").Value = 1.523408958654902
$ws.Range("N11").Value = 0.9121212171133624
$ws.Range("D12").Value = 0.007258226268384626
$ws.Range("E12").Value = 0.6141808028322584
$ws.Range("F12").Value = 1.759317199383736
$ws.Range("G12").Value = 3.919950890308719
$ws.Range("H12").Value = 7.996616307999211
$ws.Range("I12").Value = 15.68843422123317
$ws.Range("J12").Value = 52.1257047758435
$ws.Range("K12").Value = 169.5498042097182
$ws.Range("L12").Value = 547.9650619385075
$ws.Range("M12").Value = 1767.460105735589
$ws.Range("N12").Value = 1745.940780484399
$ws.Range("B13").Value = 0.08166830355129599
$ws.Range("C13").Value = 4.202182775706744
$ws.Range("D13").Value = 9.50320587664562
$ws.Range("E13").Value = 20.45561881987197
$ws.Range("F13").Value = 41.70176954636128
$ws.Range("G13").Value = 81.78882747270453
$ws.Range("H13").Value = 157.4247501335778
$ws.Range("I13").Value = 300.1339708320161
$ws.Range("J13").Value = 975.3318418770497
$ws.Range("K13").Value = 3152.302592612967
$ws.Range("L13").Value = 10167.88927098196
$ws.Range("M13").Value = 24354.03567107371
$ws.Range("N13").Value = 25460.46915666178

# --- Sheet 2: CAP_NEW ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 0.006385367863609007
$ws.Range("C2").Value = 0.1220369550162755
$ws.Range("D2").Value = 0.124902959092586
$ws.Range("E2").Value = 0.2356658493471741
$ws.Range("F2").Value = 0.4446523361176445
$ws.Range("G2").Value = 0.845351643896295
$ws.Range("H2").Value = 1.704991445856794
$ws.Range("I2").Value = 3.11160813623196
$ws.Range("J2").Value = 7.29166038338624
$ws.Range("K2").Value = 23.72049091346406
$ws.Range("L2").Value = 78.97640463456391
$ws.Range("M2").Value = 163.7101477049405
$ws.Range("N2").Value = 0.08376657732652526
$ws.Range("E3").Value = 0.03902215272332507
$ws.Range("F3").Value = 0.1898727477477608
$ws.Range("G3").Value = 0.358250298399166
$ws.Range("H3").Value = 0.6759436402826541
$ws.Range("I3").Value = 1.275364757210847
$ws.Range("J3").Value = 3.090674452706084
$ws.Range("K3").Value = 10.13985150300465
$ws.Range("L3").Value = 33.7542586744156
$ws.Range("M3").Value = 107.8133286208086
$ws.Range("N3").Value = 12.82597447194583
$ws.Range("C4").Value = 0.001664027448680314
$ws.Range("D4").Value = 0.01812410870692616
$ws.Range("E4").Value = 0.03419641538606107
$ws.Range("F4").Value = 0.06452150801816522
$ws.Range("G4").Value = 0.1217386369284503
$ws.Range("H4").Value = 0.2313594604355313
$ws.Range("I4").Value = 0.451511513611301
$ws.Range("J4").Value = 1.058843714293814
$ws.Range("K4").Value = 3.440024914077114
$ws.Range("L4").Value = 11.46443037825295
$ws.Range("M4").Value = 10.53313571320021
$ws.Range("N4").Value = 25.79960979286619
$ws.Range("C5").Value = 0.03411614863502015
$ws.Range("D5").Value = 0.05875991450236628
$ws.Range("E5").Value = 0.1108677108962861
$ws.Range("F5").Value = 0.209184261472802
$ws.Range("G5").Value = 0.394687099554677
$ws.Range("H5").Value = 0.7788084379023822
$ws.Range("I5").Value = 1.463839042550835
$ws.Range("J5").Value = 3.431536037900564
$ws.Range("K5").Value = 11.15576258712685
$ws.Range("L5").Value = 37.16141177436501
$ws.Range("M5").Value = 94.64746561101961
$ws.Range("N5").Value = 34.48099318285622
$ws.Range("C6").Value = 0.09854874261732799
$ws.Range("D6").Value = 0.1374080378394981
$ws.Range("E6").Value = 0.2592603264152458
$ws.Range("F6").Value = 0.4891701963712818
$ws.Range("G6").Value = 0.9229621991398033
$ws.Range("H6").Value = 1.839986032460754
$ws.Range("I6").Value = 3.423137222938942
$ws.Range("J6").Value = 8.023662338646133
$ws.Range("K6").Value = 26.08926650619458
$ws.Range("L6").Value = 86.89594168379152
$ws.Range("M6").Value = 171.0792121372329
$ws.Range("N6").Value = 0.5328430020948771
$ws.Range("B7").Value = 0.004888607537544012
$ws.Range("C7").Value = 0.1535568252959821
$ws.Range("D7").Value = 0.1777716031670491
$ws.Range("E7").Value = 0.3354179609076866
$ws.Range("F7").Value = 0.6328637785515783
$ws.Range("G7").Value = 1.198970645120531
$ws.Range("H7").Value = 2.40654073874037
$ws.Range("I7").Value = 4.428682641502151
$ws.Range("J7").Value = 10.37918158852545
$ws.Range("K7").Value = 33.75733421936251
$ws.Range("L7").Value = 112.4125564457501
$ws.Range("M7").Value = 283.0528681643027
$ws.Range("N7").Value = 220.154445607015
$ws.Range("B8").Value = 0.00505968530910618
$ws.Range("C8").Value = 0.4141801954178035
$ws.Range("D8").Value = 0.5327252756759985
$ws.Range("E8").Value = 1.005141555275961
$ws.Range("F8").Value = 1.896492605612812
$ws.Range("G8").Value = 3.583345931923001
$ws.Range("H8").Value = 7.165660493738692
$ws.Range("I8").Value = 13.27136133693284
$ws.Range("J8").Value = 31.1056906218102
$ws.Range("K8").Value = 101.1520333733247
$ws.Range("L8").Value = 336.881315850851
$ws.Range("M8").Value = 553.0560996026129
$ws.Range("N8").Value = 210.016913114034
$ws.Range("D9").Value = 0.009061075949674022
$ws.Range("E9").Value = 0.04952610238075569
$ws.Range("F9").Value = 0.09344543209552107
$ws.Range("G9").Value = 0.1763120528319958
$ws.Range("H9").Value = 0.3326640936504641
$ws.Range("I9").Value = 0.6367288809960798
$ws.Range("J9").Value = 1.534410664937003
$ws.Range("K9").Value = 4.984444617260823
$ws.Range("L9").Value = 16.60439790710539
$ws.Range("M9").Value = 53.00296481963424
$ws.Range("N9").Value = 73.3529984604649
$ws.Range("K10").Value = 0.06438220879409279
$ws.Range("L10").Value = 0.2074804937949305
$ws.Range("N10").Value = 0.03975999222636958
$ws.Range("K11").Value = 0.04092284375570817
$ws.Range("L11").Value = 0.1318794739876362
$ws.Range("N11").Value = 0.02527238471751815
$ws.Range("D12").Value = 0.001451645253676925
$ws.Range("E12").Value = 0.1213845153127748
$ws.Range("F12").Value = 0.2290272793102955
$ws.Range("G12").Value = 0.4321267381849967
$ws.Range("H12").Value = 0.8153330835380983
$ws.Range("I12").Value = 1.53981522790047
$ws.Range("J12").Value = 3.761676132944995
$ws.Range("K12").Value = 12.21956340319206
$ws.Range("L12").Value = 40.69610472418618
$ws.Range("M12").Value = 129.9401241477767
$ws.Range("N12").Value = 24.3059015385701
$ws.Range("B13").Value = 0.0163336607102592
$ws.Range("C13").Value = 0.8241028944310896
$ws.Range("D13").Value = 1.060204620187775
$ws.Range("E13").Value = 2.19048258864527
$ws.Range("F13").Value = 4.249230145297861
$ws.Range("G13").Value = 8.033745245978915
$ws.Range("H13").Value = 15.95128742660574
$ws.Range("I13").Value = 29.60204875987543
$ws.Range("J13").Value = 69.67733593515047
$ws.Range("K13").Value = 226.7640770895571
$ws.Range("L13").Value = 755.1861820410643
$ws.Range("M13").Value = 1566.835346521528
$ws.Range("N13").Value = 601.6184781241175

# --- Sheet 3: INVESTMENT ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 19.06517595244923
$ws.Range("C2").Value = 311.2552537690108
$ws.Range("D2").Value = 280.8592918747708
$ws.Range("E2").Value = 458.7824922166112
$ws.Range("F2").Value = 779.5022243543979
$ws.Range("G2").Value = 1318.215992942565
$ws.Range("H2").Value = 2471.333951026047
$ws.Range("I2").Value = 4168.216911052246
$ws.Range("J2").Value = 11453.50744761984
$ws.Range("K2").Value = 40038.76543247253
$ws.Range("L2").Value = 131537.0925193689
$ws.Range("M2").Value = 265412.2091320863
$ws.Range("N2").Value = 97.1717426960891
$ws.Range("E3").Value = 73.32184452407334
$ws.Range("F3").Value = 325.4058138175898
$ws.Range("G3").Value = 554.8043246158685
$ws.Range("H3").Value = 979.4017781511488
$ws.Range("I3").Value = 1720.747637724019
$ws.Range("J3").Value = 3850.733114115564
$ws.Range("K3").Value = 12208.17841258754
$ws.Range("L3").Value = 40120.64940299713
$ws.Range("M3").Value = 175590.1776582799
$ws.Range("N3").Value = 17936.72272210229
$ws.Range("C4").Value = 4.477748101928343
$ws.Range("D4").Value = 42.51807157992635
$ws.Range("E4").Value = 68.42600129504662
$ws.Range("F4").Value = 115.1360501980951
$ws.Range("G4").Value = 190.8800957719636
$ws.Range("H4").Value = 335.4480816854768
$ws.Range("I4").Value = 601.3410942880751
$ws.Range("J4").Value = 1305.15193911284
$ws.Range("K4").Value = 5819.971750632224
$ws.Range("L4").Value = 16654.64725778939
$ws.Range("M4").Value = 17021.65264388867
$ws.Range("N4").Value = 41264.06386169812
$ws.Range("C5").Value = 91.80348552346207
$ws.Range("D5").Value = 137.8472338276812
$ws.Range("E5").Value = 221.8429634721417
$ws.Range("F5").Value = 373.2809472277563
$ws.Range("G5").Value = 618.8496377467558
$ws.Range("H5").Value = 1129.194354114664
$ws.Range("I5").Value = 1949.599390430903
$ws.Range("J5").Value = 4229.779951036992
$ws.Range("K5").Value = 13501.26167107028
$ws.Range("L5").Value = 50373.4987488886
$ws.Range("M5").Value = 152192.6528157138
$ws.Range("N5").Value = 53332.95291666834
$ws.Range("C6").Value = 251.348568045495
$ws.Range("D6").Value = 308.9784620466522
$ws.Range("E6").Value = 504.7150404488798
$ws.Range("F6").Value = 857.5447044506393
$ws.Range("G6").Value = 1439.239564472635
$ws.Range("H6").Value = 2667.004554470889
$ws.Range("I6").Value = 4585.531929732318
$ws.Range("J6").Value = 9937.546516283395
$ws.Range("K6").Value = 41268.56137399755
$ws.Range("L6").Value = 145472.4959728354
$ws.Range("M6").Value = 276055.7877367791
$ws.Range("N6").Value = 618.1138677201203
$ws.Range("B7").Value = 16.39140330123432
$ws.Range("C7").Value = 429.378666029131
$ws.Range("D7").Value = 430.01706404887
$ws.Range("E7").Value = 684.7993715279603
$ws.Range("F7").Value = 1144.224040259039
$ws.Range("G7").Value = 1887.647393971313
$ws.Range("H7").Value = 3490.013510136059
$ws.Range("I7").Value = 5872.654616763927
$ws.Range("J7").Value = 12747.50324339518
$ws.Range("K7").Value = 57210.5797416068
$ws.Range("L7").Value = 188854.0061290132
$ws.Range("M7").Value = 455298.0935393163
$ws.Range("N7").Value = 354844.9354293869
$ws.Range("B8").Value = 17.86063854429172
$ws.Range("C8").Value = 1131.262793150509
$ws.Range("D8").Value = 1243.860246175889
$ws.Range("E8").Value = 1948.436750655791
$ws.Range("F8").Value = 3303.121171195834
$ws.Range("G8").Value = 5535.982797146482
$ws.Range("H8").Value = 10370.50214956332
$ws.Range("I8").Value = 17910.76383309782
$ws.Range("J8").Value = 38963.29912978567
$ws.Range("K8").Value = 170316.7592330028
$ws.Range("L8").Value = 461501.979952524
$ws.Range("M8").Value = 688909.3768437783
$ws.Range("N8").Value = 311817.7787036936
$ws.Range("D9").Value = 23.68166565903003
$ws.Range("E9").Value = 106.4840916847676
$ws.Range("F9").Value = 174.8177143643008
$ws.Range("G9").Value = 280.6076845642347
$ws.Range("H9").Value = 482.7188863733789
$ws.Range("I9").Value = 834.4968714334622
$ws.Range("J9").Value = 1866.395756402773
$ws.Range("K9").Value = 6070.704632700475
$ws.Range("L9").Value = 20330.75688541798
$ws.Range("M9").Value = 64709.89632939233
$ws.Range("N9").Value = 83986.2491172939
$ws.Range("K10").Value = 108.4846656401343
$ws.Range("L10").Value = 345.7849161536932
$ws.Range("N10").Value = 64.64338576132072
$ws.Range("K11").Value = 68.95540095680582
$ws.Range("L11").Value = 219.7890125530545
$ws.Range("N11").Value = 41.08885396912972
$ws.Range("D12").Value = 3.723324911155946
$ws.Range("E12").Value = 257.6932567832552
$ws.Range("F12").Value = 424.8685058485293
$ws.Range("G12").Value = 685.893165184136
$ws.Range("H12").Value = 1182.934157582085
$ws.Range("I12").Value = 2024.025524466052
$ws.Range("J12").Value = 6411.776968604743
$ws.Range("K12").Value = 20780.71171910246
$ws.Range("L12").Value = 69170.37165583186
$ws.Range("M12").Value = 208342.0968548207
$ws.Range("N12").Value = 38945.08064702537
$ws.Range("B13").Value = 53.31721779797528
$ws.Range("C13").Value = 2219.526514619537
$ws.Range("D13").Value = 2471.485360123976
$ws.Range("E13").Value = 4324.501812608527
$ws.Range("F13").Value = 7497.901171716181
$ws.Range("G13").Value = 12512.12065641595
$ws.Range("H13").Value = 23108.55142310307
$ws.Range("I13").Value = 39667.37780898883
$ws.Range("J13").Value = 90765.694066357
$ws.Range("K13").Value = 367392.9340337696
$ws.Range("L13").Value = 1124581.072453373
$ws.Range("M13").Value = 2303531.943554055
$ws.Range("N13").Value = 902948.8012480151

# --- Sheet 4: REMOVAL ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 0.02915693088405939
$ws.Range("C2").Value = 0.5864033008213906
$ws.Range("D2").Value = 1.156736447362879
$ws.Range("E2").Value = 2.232836216071437
$ws.Range("F2").Value = 4.263212180078946
$ws.Range("G2").Value = 8.094108417671119
$ws.Range("H2").Value = 15.32221111557303
$ws.Range("I2").Value = 28.96013429885784
$ws.Range("J2").Value = 93.45935233249475
$ws.Range("K2").Value = 301.3169244422875
$ws.Range("L2").Value = 971.166445111608
$ws.Range("M2").Value = 2324.628290376588
$ws.Range("N2").Value = 1856.45768584063
$ws.Range("E3").Value = 0.1781833457686076
$ws.Range("F3").Value = 1.045182193931899
$ws.Range("G3").Value = 2.681028305343615
$ws.Range("H3").Value = 5.767528945903679
$ws.Range("I3").Value = 11.59111231216325
$ws.Range("J3").Value = 39.20476700994847
$ws.Range("K3").Value = 128.1935547432587
$ws.Range("L3").Value = 414.9721224088885
$ws.Range("M3").Value = 1339.15537028974
$ws.Range("N3").Value = 1255.857830410614
$ws.Range("C4").Value = 0.007598298852421526
$ws.Range("D4").Value = 0.09035678609865967
$ws.Range("E4").Value = 0.2465048015601257
$ws.Range("F4").Value = 0.5411235596339394
$ws.Range("G4").Value = 1.09700774652184
$ws.Range("H4").Value = 2.145845340069105
$ws.Range("I4").Value = 4.124783262006889
$ws.Range("J4").Value = 13.49113147749768
$ws.Range("K4").Value = 43.67547260534667
$ws.Range("L4").Value = 140.948649766913
$ws.Range("M4").Value = 216.5988908538056
$ws.Range("N4").Value = 384.1548921935422
$ws.Range("C5").Value = 0.1557815006165304
$ws.Range("D5").Value = 0.4240916124994815
$ws.Range("E5").Value = 0.9303368677336648
$ws.Range("F5").Value = 1.885516143865181
$ws.Range("G5").Value = 3.687740342744984
$ws.Range("H5").Value = 7.08816175492472
$ws.Range("I5").Value = 13.50404818436978
$ws.Range("J5").Value = 43.85844194769599
$ws.Range("K5").Value = 141.6796370390724
$ws.Range("L5").Value = 456.9218482232159
$ws.Range("M5").Value = 1254.673595241534
$ws.Range("N5").Value = 1348.942143206012
$ws.Range("C6").Value = 0.4499942585266118
$ws.Range("D6").Value = 1.077428221264046
$ws.Range("E6").Value = 2.261265328182977
$ws.Range("F6").Value = 4.494919192892025
$ws.Range("G6").Value = 8.709358458370579
$ws.Range("H6").Value = 16.66112690514422
$ws.Range("I6").Value = 31.66445651747043
$ws.Range("J6").Value = 102.6392476256503
$ws.Range("K6").Value = 331.36524168425
$ws.Range("L6").Value = 1068.466112870403
$ws.Range("M6").Value = 2475.066548622116
$ws.Range("N6").Value = 1964.01850211163
$ws.Range("B7").Value = 0.02232240884723293
$ws.Range("C7").Value = 0.7234951270937264
$ws.Range("D7").Value = 1.535237607308562
$ws.Range("E7").Value = 3.066826469900738
$ws.Range("F7").Value = 5.956615413058629
$ws.Range("G7").Value = 11.40904480841473
$ws.Range("H7").Value = 21.69664258669961
$ws.Range("I7").Value = 41.10719527316126
$ws.Range("J7").Value = 132.9177584095422
$ws.Range("K7").Value = 428.7898841399463
$ws.Range("L7").Value = 1382.278369492012
$ws.Range("M7").Value = 3765.699467760128
$ws.Range("N7").Value = 5108.799926887605
$ws.Range("B8").Value = 0.02310358588632959
$ws.Range("C8").Value = 1.914337354917396
$ws.Range("D8").Value = 4.346872860287252
$ws.Range("E8").Value = 8.936560327300771
$ws.Range("F8").Value = 17.59634391457388
$ws.Range("G8").Value = 33.93555052011678
$ws.Range("H8").Value = 64.76422768139936
$ws.Range("I8").Value = 122.9315156323438
$ws.Range("J8").Value = 398.082170360828
$ws.Range("K8").Value = 1284.792874767451
$ws.Range("L8").Value = 4142.339676119458
$ws.Range("M8").Value = 8589.163763836768
$ws.Range("N8").Value = 8506.974161114817
$ws.Range("D9").Value = 0.04137477602590878
$ws.Range("E9").Value = 0.2675213622394051
$ws.Range("F9").Value = 0.694212832995209
$ws.Range("G9").Value = 1.499290699807975
$ws.Range("H9").Value = 3.018304826065802
$ws.Range("I9").Value = 5.88436786280738
$ws.Range("J9").Value = 19.4577583260288
$ws.Range("K9").Value = 63.19986982747471
$ws.Range("L9").Value = 204.1648206103627
$ws.Range("M9").Value = 658.4436989531511
$ws.Range("N9").Value = 1229.754906243396
$ws.Range("K10").Value = 0.5879653771150026
$ws.Range("L10").Value = 2.482764407205692
$ws.Range("M10").Value = 2.188781718648191
$ws.Range("N10").Value = 1.310504466884336
$ws.Range("K11").Value = 0.3737246005087504
$ws.Range("L11").Value = 1.578103358386706
$ws.Range("M11").Value = 1.391241058132331
$ws.Range("N11").Value = 0.8329874128889153
$ws.Range("D12").Value = 0.006628517140077284
$ws.Range("E12").Value = 0.5608957103490945
$ws.Range("F12").Value = 1.606682373866426
$ws.Range("G12").Value = 3.579863826765953
$ws.Range("H12").Value = 7.302845943378275
$ws.Range("I12").Value = 14.3273371883408
$ws.Range("J12").Value = 47.60338335693471
$ws.Range("K12").Value = 154.8400038444915
$ws.Range("L12").Value = 500.4247140990935
$ws.Range("M12").Value = 1614.118818023369
$ws.Range("N12").Value = 1594.466466195798
$ws.Range("B13").Value = 0.07458292561762192
$ws.Range("C13").Value = 3.837609840828077
$ws.Range("D13").Value = 8.678726827986868
$ws.Range("E13").Value = 18.68093042910682
$ws.Range("F13").Value = 38.08380780489613
$ws.Range("G13").Value = 74.69299312575758
$ws.Range("H13").Value = 143.7668950991578
$ws.Range("I13").Value = 274.0949505315215
$ws.Range("J13").Value = 890.7140108466208
$ws.Range("K13").Value = 2878.815153071203
$ws.Range("L13").Value = 9285.743626467547
$ws.Range("M13").Value = 22241.12846673398
$ws.Range("N13").Value = 23251.57000608382
